# Apply the latest cryptos.xlsx data refresh (GitHub Actions bot update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds prices formatted as plain text (e.g. "1.00", "64.820.44").
# Force text format on cells whose new value would otherwise be auto-parsed as
# a number by Excel (losing trailing zeros / becoming numeric), so the displayed
# text matches the source feed exactly.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '64.820.44'
$ws.Range("E2").Value = '  +1.54%  '

# Row 3
$ws.Range("D3").Value = '2.633.34'
$ws.Range("E3").Value = '  +0.38%  '

# Row 4
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.19%  '

# Row 5
$ws.Range("D5").Value = '595.44'
$ws.Range("E5").Value = '  -0.58%  '

# Row 6
$ws.Range("D6").Value = '154.51'
$ws.Range("E6").Value = '  +1.93%  '

# Row 7
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.05%  '

# Row 8
$ws.Range("E8").Value = '  +0.25%  '

# Row 9
$ws.Range("E9").Value = '  +5.04%  '

# Row 10
$ws.Range("D10").Value = '0.399'
$ws.Range("E10").Value = '  +3.39%  '

# Row 11
$ws.Range("E11").Value = '  +1.77%  '

# Row 12
$ws.Range("E12").Value = '  +1.24%  '

# Row 13
$ws.Range("D13").Value = '28.83'
$ws.Range("E13").Value = '  +3.28%  '

# Row 14
$ws.Range("D14").Value = '3.105.90'
$ws.Range("E14").Value = '  +0.46%  '

# Row 15
$ws.Range("D15").Value = '64.703.22'
$ws.Range("E15").Value = '  +1.66%  '

# Row 16
$ws.Range("E16").Value = '  +12.11%  '

# Row 17
$ws.Range("D17").Value = '2.631.93'
$ws.Range("E17").Value = '  +0.92%  '

# Row 18
$ws.Range("D18").Value = '12.49'
$ws.Range("E18").Value = '  +0.76%  '

# Row 19
$ws.Range("E19").Value = '  +1.94%  '

# Row 20
$ws.Range("D20").Value = '354.12'
$ws.Range("E20").Value = '  +1.43%  '

# Row 21
$ws.Range("D21").Value = '7.20'
$ws.Range("E21").Value = '  +4.16%  '

# Row 22
$ws.Range("E22").Value = '  +0.36%  '

# Row 23
$ws.Range("D23").Value = '67.73'
$ws.Range("E23").Value = '  +0.90%  '

# Row 24
$ws.Range("E24").Value = '  -0.15%  '

# Row 25
$ws.Range("D25").Value = '9.46'
$ws.Range("E25").Value = '  +0.54%  '

# Row 26
$ws.Range("E26").Value = '  -1.59%  '

# Row 27
$ws.Range("E27").Value = '  +2.60%  '

# Row 28
$ws.Range("E28").Value = '  +1.63%  '

# Row 29
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  +0.17%  '

# Row 30
$ws.Range("D30").Value = '530.99'
$ws.Range("E30").Value = '  -4.89%  '

# Row 31
$ws.Range("D31").Value = '0.0₃0915'
$ws.Range("E31").Value = '  +7.31%  '

# Row 32
$ws.Range("E32").Value = '  -0.39%  '

# Row 33
$ws.Range("E33").Value = '  +4.37%  '

# Row 34
$ws.Range("E34").Value = '  +8.94%  '

# Row 35
$ws.Range("E35").Value = '  +0.97%  '

# Row 36
$ws.Range("D36").Value = '0.426'
$ws.Range("E36").Value = '  +1.71%  '

# Row 37
$ws.Range("D37").Value = '164.83'
$ws.Range("E37").Value = '  -1.61%  '

# Row 38
$ws.Range("B38").Value = 'EthereumClassic'
$ws.Range("C38").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D38").Value = '20.18'
$ws.Range("E38").Value = '  +2.80%  '

# Row 39
$ws.Range("B39").Value = 'Stacks'
$ws.Range("C39").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D39").Value = '2.01'
$ws.Range("E39").Value = '  +3.30%  '

# Row 40
$ws.Range("D40").Value = '1.00'
$ws.Range("E40").Value = '  +0.29%  '

# Row 41
$ws.Range("E41").Value = '  +0.00%  '

# Row 42
$ws.Range("D42").Value = '166.67'
$ws.Range("E42").Value = '  -0.89%  '

# Row 43
$ws.Range("D43").Value = '41.93'

# Row 44
$ws.Range("D44").Value = '4.09'
$ws.Range("E44").Value = '  +2.85%  '

# Row 45
$ws.Range("D45").Value = '23.37'
$ws.Range("E45").Value = '  +6.39%  '

# Row 46
$ws.Range("D46").Value = '0.0601'
$ws.Range("E46").Value = '  +1.68%  '

# Row 47
$ws.Range("D47").Value = '2.22'
$ws.Range("E47").Value = '  +9.09%  '

# Row 48
$ws.Range("E48").Value = '  +1.55%  '

# Row 49
$ws.Range("D49").Value = '0.0251'
$ws.Range("E49").Value = '  -0.73%  '

# Row 50
$ws.Range("E50").Value = '  +1.50%  '

# Row 51
$ws.Range("D51").Value = '19.42'
$ws.Range("E51").Value = '  +0.00%  '
